$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan breakdown")

$ws.Range("A10").Value = "LookupValue"
$ws.Range("B10").Value = "Permissions"
$ws.Range("A11").Value = "LookupValue"
$ws.Range("B11").Value = "Copy Document"

$table = $ws.ListObjects.Item("Table29")
$table.Resize($ws.Range("A1:B11"))

$ws.Range("C16").Select()
